$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2321981424148607
$ws.Range("C2").Value = 0.4922600619195047
$ws.Range("J2").Value = 0.01857585139318885
$ws.Range("P2").Value = 0.151702786377709
$ws.Range("S2").Value = 0.1052631578947368
$ws.Range("B3").Value = 0.005988023952095809
$ws.Range("C3").Value = 0.03592814371257485
$ws.Range("J3").Value = 0.05389221556886228
$ws.Range("P3").Value = 0.7604790419161677
$ws.Range("S3").Value = 0.1437125748502994
$ws.Range("J4").Value = 0.04761904761904762
$ws.Range("P4").Value = 0.6190476190476191
$ws.Range("S4").Value = 0.3333333333333333
$ws.Range("B6").Value = 0.05882352941176471
$ws.Range("D6").Value = 0.0053475935828877
$ws.Range("F6").Value = 0.0213903743315508
$ws.Range("J6").Value = 0.2887700534759358
$ws.Range("O6").Value = 0.0160427807486631
$ws.Range("Q6").Value = 0.1871657754010695
$ws.Range("R6").Value = 0.0855614973262032
$ws.Range("S6").Value = 0.3368983957219251
$ws.Range("B7").Value = 0.103960396039604
$ws.Range("F7").Value = 0.0396039603960396
$ws.Range("J7").Value = 0.1782178217821782
$ws.Range("O7").Value = 0.0198019801980198
$ws.Range("Q7").Value = 0.1831683168316832
$ws.Range("R7").Value = 0.1237623762376238
$ws.Range("S7").Value = 0.3514851485148515
$ws.Range("B8").Value = 0.0979020979020979
$ws.Range("D8").Value = 0.01398601398601399
$ws.Range("E8").Value = 0.002331002331002331
$ws.Range("F8").Value = 0.04662004662004662
$ws.Range("J8").Value = 0.1608391608391608
$ws.Range("O8").Value = 0.01864801864801865
$ws.Range("Q8").Value = 0.1818181818181818
$ws.Range("R8").Value = 0.09324009324009325
$ws.Range("S8").Value = 0.3846153846153846
$ws.Range("B9").Value = 0.06315789473684211
$ws.Range("D9").Value = 0.01578947368421053
$ws.Range("F9").Value = 0.01052631578947368
$ws.Range("J9").Value = 0.2157894736842105
$ws.Range("O9").Value = 0.01052631578947368
$ws.Range("Q9").Value = 0.1894736842105263
$ws.Range("R9").Value = 0.08947368421052632
$ws.Range("S9").Value = 0.4052631578947368
$ws.Range("B10").Value = 0.1026814911706998
$ws.Range("D10").Value = 0.02027468933943754
$ws.Range("E10").Value = 0.0006540222367560497
$ws.Range("F10").Value = 0.06147809025506867
$ws.Range("J10").Value = 0.1530412034009156
$ws.Range("O10").Value = 0.01177240026160889
$ws.Range("Q10").Value = 0.2452583387835186
$ws.Range("R10").Value = 0.07521255722694571
$ws.Range("S10").Value = 0.3296272073250491
$ws.Range("G11").Value = 0.1755485893416928
$ws.Range("J11").Value = 0.09404388714733543
$ws.Range("K11").Value = 0.2225705329153605
$ws.Range("L11").Value = 0.5078369905956113
$ws.Range("F12").Value = 0.006134969325153374
$ws.Range("G12").Value = 0.7300613496932515
$ws.Range("J12").Value = 0.196319018404908
$ws.Range("K12").Value = 0.01840490797546012
$ws.Range("L12").Value = 0.03680981595092025
$ws.Range("S12").Value = 0.01226993865030675
$ws.Range("G13").Value = 0.6153846153846154
$ws.Range("J13").Value = 0.3269230769230769
$ws.Range("S13").Value = 0.0576923076923077
$ws.Range("F15").Value = 0.01428571428571429
$ws.Range("H15").Value = 0.1714285714285714
$ws.Range("I15").Value = 0.05714285714285714
$ws.Range("J15").Value = 0.3666666666666666
$ws.Range("K15").Value = 0.06666666666666667
$ws.Range("M15").Value = 0.03333333333333333
$ws.Range("O15").Value = 0.0380952380952381
$ws.Range("S15").Value = 0.2523809523809524
$ws.Range("F16").Value = 0.005076142131979695
$ws.Range("H16").Value = 0.116751269035533
$ws.Range("I16").Value = 0.06091370558375635
$ws.Range("J16").Value = 0.4619289340101523
$ws.Range("K16").Value = 0.09137055837563451
$ws.Range("M16").Value = 0.04060913705583756
$ws.Range("O16").Value = 0.05076142131979695
$ws.Range("S16").Value = 0.1725888324873096
$ws.Range("F17").Value = 0.01612903225806452
$ws.Range("H17").Value = 0.1577060931899641
$ws.Range("I17").Value = 0.09498207885304659
$ws.Range("J17").Value = 0.4193548387096774
$ws.Range("K17").Value = 0.09498207885304659
$ws.Range("M17").Value = 0.01971326164874552
$ws.Range("O17").Value = 0.07168458781362007
$ws.Range("S17").Value = 0.1254480286738351
$ws.Range("F18").Value = 0.009615384615384616
$ws.Range("H18").Value = 0.1586538461538461
$ws.Range("I18").Value = 0.09134615384615384
$ws.Range("J18").Value = 0.4423076923076923
$ws.Range("K18").Value = 0.07211538461538461
$ws.Range("M18").Value = 0.009615384615384616
$ws.Range("N18").Value = 0.004807692307692308
$ws.Range("O18").Value = 0.05288461538461538
$ws.Range("S18").Value = 0.1586538461538461
$ws.Range("F19").Value = 0.01001540832049307
$ws.Range("H19").Value = 0.1926040061633282
$ws.Range("I19").Value = 0.07318952234206472
$ws.Range("J19").Value = 0.3975346687211094
$ws.Range("K19").Value = 0.1140215716486903
$ws.Range("M19").Value = 0.02080123266563945
$ws.Range("O19").Value = 0.06240369799691833
$ws.Range("S19").Value = 0.1294298921417565
